# Applies the "cryptos list" data refresh described by the commit:
# "Updated cryptos list on Fri Sep 29 14:57:14 UTC 2023 with GitHub Actions"
#
# Updates Price (col D) and Volume(1h) (col E) values for the crypto rows,
# including one row reorder where the Toncoin/Avalanche rows swap places
# (rows 23 and 24 get each other's Coin/Link/Price/Volume data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting before writing so numeric-looking strings (e.g. "1.00",
# "65.79") are preserved verbatim as text instead of being coerced to numbers,
# matching the original inlineStr/text cell type used throughout column D & E.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.953.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.675.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.02%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.696.79"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.79"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.019.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.86"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.462.60"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.13%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.97"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.970"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.821.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.71"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.11%  "
